$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: "List of all entries" / "additional shared assetts in listview_userentries.xml" ---
$ws.Range("B16").Value = 14
$ws.Range("C16").Value = "List of all entries"
$ws.Range("C16").WrapText = $true

$ws.Range("D16").Value = "additional shared assetts in listview_userentries.xml"
$ws.Range("D16").WrapText = $true
$ws.Range("D16").Font.Name = "Bitstream Vera Sans Mono"
$ws.Range("D16").Font.Family = 1
$ws.Range("D16").Font.Size = 9
$ws.Range("D16").Font.Color = 8355711

# --- Row 17: "Hide closed items" ---
$ws.Range("B17").Value = 15
$ws.Range("C17").Value = "Hide closed items"
$ws.Range("C17").WrapText = $true

# Leave the cursor one row below the last entered row, like Excel does
# after typing values down a column.
$ws.Range("C18").Select()
